# update excel upload api
# Adds a new "Image" column (O) to the artist roster sheet, giving each
# artist row the filename(s) of their uploaded image(s), resizes a few
# rows that now wrap onto more lines, adds the narrow column P width
# override that shipped alongside this change, and moves the saved
# window selection/scroll position to reflect where the author was
# working (A28 / I26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column O header -------------------------------------------------
$ws.Cells.Item(1, 15).Value = "Image"

# --- Column O data: artist image filename(s) per row -----------------
$images = @{
    2  = "artist-1.jpg"
    3  = "artist-2.jpg"
    4  = "artist-3.jpg"
    5  = "artist-4.jpg"
    6  = "artist-5.jpg,artist-62.jpg,artist-61.jpg"
    7  = "artist-6.jpg,artist-60.jpg"
    8  = "artist-7.jpg"
    9  = "artist-8.jpg"
    10 = "artist-9.jpg"
    11 = "artist-10.jpg"
    12 = "artist-11.jpg"
    13 = "artist-12.jpg"
    14 = "artist-13.jpg"
    15 = "artist-14.jpg,artist-58.jpg"
    16 = "artist-15.jpg"
    17 = "artist-16.jpg,artist-59.jpg"
    18 = "artist-17.jpg"
    19 = "artist-18.jpg"
    20 = "artist-19.jpg"
    21 = "artist-20.jpg"
    22 = "artist-21.jpg"
    23 = "artist-22.jpg"
    24 = "artist-23.jpg"
    25 = "artist-24.jpg"
    26 = "artist-25.jpg,artist-57.jpg"
    27 = "artist-26.jpg"
    28 = "artist-27.jpg"
}

foreach ($row in 2..28) {
    $cell = $ws.Cells.Item($row, 15)
    $cell.Value = $images[$row]
    $cell.WrapText = $true
}

# --- Row heights that changed because of the new wrapped content -----
$ws.Rows.Item(15).RowHeight = 48
$ws.Rows.Item(17).RowHeight = 48
$ws.Rows.Item(19).RowHeight = 32
$ws.Rows.Item(25).RowHeight = 32

# --- New narrow column width override (column P) ---------------------
$ws.Columns.Item(16).ColumnWidth = 9.92

# --- Window scroll position / selection -------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 26
$ws.Range("A28").Select()
